$wb = $excel.ActiveWorkbook

# --- TestCases sheet edits -------------------------------------------------
$ws = $wb.Worksheets.Item("TestCases")

# Row 9: Grouping column changes from "Sanity" to "Regression"
$ws.Range("D9").Value = "Regression"

# Row 13 (AddressPage): TestClass/TestCase columns are cleared out
$ws.Range("B13").ClearContents()
$ws.Range("C13").ClearContents()

# Row 15 (PaymentPage): now carries the end-to-end test info that used to
# live in the rows (16/17) that get removed below
$ws.Range("B15").Value = "EndToEndTest"
$ws.Range("C15").Value = "endToEndTest"
$ws.Range("D15").Value = "Regression"
$ws.Range("E15").Value = "User should be able to order the product"

# Rows 16 (OrderSummary) and 17 (OrderConfirmationPage) are removed entirely
$ws.Rows("16:17").Delete()

# Make TestCases the active sheet/tab (was ProductDetails before)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select()

# --- ProductDetails sheet: no longer the selected tab ----------------------
# (handled automatically by activating TestCases above)
